$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.556.38"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.789.82"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +4.15%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.50"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.77%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9989"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.05%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5359"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +9.62%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3771"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +7.93%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.14"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.73%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07484"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.99%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.107"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +5.86%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9986"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.06%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.88"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +4.96%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.154"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +5.11%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.787.42"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.18%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.065"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "90.54"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +4.48%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001062"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.31%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06467"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.32%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9987"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.15%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.91"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.36%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.936"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +5.18%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.586.13"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.24"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +4.28%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.094"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.52%  "

$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.49"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.66%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.45"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.89%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.387"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +14.72%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.989.65"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +4.20%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "121.76"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.76%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.109"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +7.55%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1026"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +11.12%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.676"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +6.17%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.622"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.69%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02277"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.36%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "8.568"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +14.52%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06015"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.33%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "11.42"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.96%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.975"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +5.17%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2082"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6227"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.06%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.412"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.61%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9983"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.12%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +4.90%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.41"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.80%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5846"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.62%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.635"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.78%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "121.71"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.22%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.910"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +4.03%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.129"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.58%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06751"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.21%  "
